# Swap the Batch/MRP/Qty/Value columns (B, E, F, G) between each pair of
# adjacent rows listed below. Columns A, C, D and H:M are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(149,150),
    @(279,280),
    @(313,314),
    @(346,347),
    @(350,351),
    @(355,356),
    @(372,373),
    @(375,376),
    @(382,383),
    @(419,420),
    @(421,422),
    @(581,582),
    @(583,584),
    @(590,591),
    @(599,600),
    @(720,721),
    @(872,873)
)

$cols = @(2, 5, 6, 7)  # B, E, F, G

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($c in $cols) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value2 = $v2
        $ws.Cells.Item($r2, $c).Value2 = $v1
    }
}
